# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price refresh values to the Sagittarius_Profits workbook.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1876.8148
$ws.Range("I40").Value = 1798.9048
$ws.Range("K40").Value = 1798.9048
$ws.Range("M40").Value = -1623.9048
$ws.Range("H41").Value = 199.66667
$ws.Range("I41").Value = 199.66667
$ws.Range("K41").Value = 199.66667
$ws.Range("M41").Value = 240.33333
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H75").Value = 28350
$ws.Range("J75").Value = 28350
$ws.Range("L75").Value = 28350
$ws.Range("N75").Value = -30222
$ws.Range("H78").Value = 28350
$ws.Range("J78").Value = 28350
$ws.Range("L78").Value = 85050
$ws.Range("N78").Value = -94410
$ws.Range("H99").Value = 547
$ws.Range("I99").Value = 285
$ws.Range("K99").Value = 855
$ws.Range("M99").Value = 643
$ws.Range("H101").Value = 33333516
$ws.Range("J101").Value = 300
$ws.Range("L101").Value = 900
$ws.Range("N101").Value = -4144
$ws.Range("H118").Value = 875.2857
$ws.Range("I118").Value = 875.2857
$ws.Range("K118").Value = 2625.8571
$ws.Range("M118").Value = -968.8571000000002
$ws.Range("H135").Value = 269.7857
$ws.Range("I135").Value = 300.16666
$ws.Range("J135").Value = 87.5
$ws.Range("K135").Value = 2701.49994
$ws.Range("L135").Value = 787.5
$ws.Range("M135").Value = -166.4999399999997
$ws.Range("N135").Value = -5857.5
$ws.Range("H137").Value = 1261.3125
$ws.Range("I137").Value = 1260.9166
$ws.Range("J137").Value = 1262.5
$ws.Range("K137").Value = 3782.7498
$ws.Range("L137").Value = 3787.5
$ws.Range("M137").Value = -1232.7498
$ws.Range("N137").Value = -8887.5
$ws.Range("H138").Value = 4865.054
$ws.Range("I138").Value = 2622.6365
$ws.Range("J138").Value = 5813.769
$ws.Range("K138").Value = 7867.9095
$ws.Range("L138").Value = 17441.307
$ws.Range("M138").Value = -2727.9095
$ws.Range("N138").Value = -27721.307

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 253333
$ws.Range("I6").Value = 253333
$ws.Range("K6").Value = 253333
$ws.Range("M6").Value = -253160
$ws.Range("H16").Value = 975
$ws.Range("I16").Value = 800
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -513
$ws.Range("N16").Value = -1574

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1022.1
$ws.Range("I94").Value = 656.8
$ws.Range("K94").Value = 656.8
$ws.Range("M94").Value = -205.8
$ws.Range("H107").Value = 2394.1667
$ws.Range("I107").Value = 2029.3334
$ws.Range("K107").Value = 2029.3334
$ws.Range("M107").Value = -109.3334

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1994.6666
$ws.Range("I31").Value = 2032.1578
$ws.Range("K31").Value = 2032.1578
$ws.Range("M31").Value = -1737.1578
$ws.Range("H34").Value = 1994.6666
$ws.Range("I34").Value = 2032.1578
$ws.Range("K34").Value = 2032.1578
$ws.Range("M34").Value = -1830.1578
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H86").Value = 11000.375
$ws.Range("I86").Value = 11999.5
$ws.Range("J86").Value = 10001.25
$ws.Range("K86").Value = 11999.5
$ws.Range("L86").Value = 10001.25
$ws.Range("M86").Value = -10876.5
$ws.Range("N86").Value = -12247.25
$ws.Range("H89").Value = 11000.375
$ws.Range("I89").Value = 11999.5
$ws.Range("J89").Value = 10001.25
$ws.Range("K89").Value = 59997.5
$ws.Range("L89").Value = 50006.25
$ws.Range("M89").Value = -54381.5
$ws.Range("N89").Value = -61238.25
$ws.Range("H94").Value = 81711.2
$ws.Range("I94").Value = 144083.38
$ws.Range("K94").Value = 144083.38
$ws.Range("M94").Value = -143632.38
$ws.Range("H129").Value = 139969
$ws.Range("J129").Value = 139969
$ws.Range("L129").Value = 139969
$ws.Range("N129").Value = -149969

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 848.2381
$ws.Range("J12").Value = 875.75
$ws.Range("L12").Value = 2627.25
$ws.Range("N12").Value = -2973.25
$ws.Range("H23").Value = 1653.2
$ws.Range("I23").Value = 1859.6666
$ws.Range("J23").Value = 1564.7142
$ws.Range("K23").Value = 5578.9998
$ws.Range("L23").Value = 4694.142599999999
$ws.Range("M23").Value = -5343.9998
$ws.Range("N23").Value = -5164.142599999999
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H33").Value = 50
$ws.Range("I33").Value = 50
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 300
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -17
$ws.Range("N33").ClearContents()
$ws.Range("H81").Value = 2333.3333
$ws.Range("J81").Value = 2500
$ws.Range("L81").Value = 7500
$ws.Range("N81").Value = -9746
$ws.Range("H84").Value = 2333.3333
$ws.Range("J84").Value = 2500
$ws.Range("L84").Value = 22500
$ws.Range("N84").Value = -33732
$ws.Range("H132").Value = 1689.5
$ws.Range("J132").Value = 1689.5
$ws.Range("L132").Value = 15205.5
$ws.Range("N132").Value = -20265.5

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 34.615383
$ws.Range("I2").Value = 43.857143
$ws.Range("J2").Value = 23.833334
$ws.Range("K2").Value = 43.857143
$ws.Range("L2").Value = 23.833334
$ws.Range("M2").Value = 69.14285699999999
$ws.Range("N2").Value = -249.833334
$ws.Range("H70").Value = 5798.375
$ws.Range("I70").Value = 5484
$ws.Range("K70").Value = 5484
$ws.Range("M70").Value = -5214
$ws.Range("H73").Value = 5798.375
$ws.Range("I73").Value = 5484
$ws.Range("K73").Value = 5484
$ws.Range("M73").Value = -4548
$ws.Range("H107").Value = 1622.7587
$ws.Range("I107").Value = 1030.9412
$ws.Range("J107").Value = 2461.1667
$ws.Range("K107").Value = 1030.9412
$ws.Range("M107").Value = 889.0588
$ws.Range("N107").Value = -6301.1667

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 64077
$ws.Range("I63").Value = 64077
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 64077
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -63328
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 64077
$ws.Range("I66").Value = 64077
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 192231
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -188487
$ws.Range("N66").ClearContents()
$ws.Range("H132").Value = 5122.933
$ws.Range("I132").Value = 5295.6924
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 15887.0772
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -13357.0772
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 1261.5714
$ws.Range("I136").Value = 1344.1818
$ws.Range("K136").Value = 4032.5454
$ws.Range("M136").Value = -1482.5454

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 825.25
$ws.Range("I23").Value = 500.33334
$ws.Range("K23").Value = 500.33334
$ws.Range("M23").Value = -271.33334
$ws.Range("H56").Value = 32499
$ws.Range("J56").Value = 32499
$ws.Range("L56").Value = 32499
$ws.Range("N56").Value = -33927
$ws.Range("H81").Value = 1429678.4
$ws.Range("I81").Value = 949.4
$ws.Range("J81").Value = 5001501
$ws.Range("K81").Value = 1898.8
$ws.Range("L81").Value = 10003002
$ws.Range("M81").Value = -837.8
$ws.Range("N81").Value = -10005124
$ws.Range("H84").Value = 1429678.4
$ws.Range("I84").Value = 949.4
$ws.Range("J84").Value = 5001501
$ws.Range("K84").Value = 9494
$ws.Range("L84").Value = 50015010
$ws.Range("M84").Value = -4190
$ws.Range("N84").Value = -50025618
$ws.Range("H110").Value = 34793.668
$ws.Range("J110").Value = 34793.668
$ws.Range("L110").Value = 34793.668
$ws.Range("N110").Value = -42973.668
